$d = $word.ActiveDocument

# --- Locate the existing "Chapter 5" heading paragraph (it will become "Chapter 4") ---
$findRng = $d.Content
$findRng.Find.Execute("Chapter 5") | Out-Null
$headingPara = $findRng.Paragraphs(1)
$pRange = $headingPara.Range
$cursor = $d.Range($pRange.End, $pRange.End)

# --- Insert the new Chapter 4 body content plus a brand-new "Chapter 5" heading paragraph.
#     A trailing placeholder paragraph is appended too, since Range.InsertXML merges the
#     last inserted paragraph''s content into whatever paragraph follows the insertion point;
#     the placeholder absorbs that merge so our real paragraphs stay intact. ---
$xmlSnippet = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">FROM: Screwtape [mailto: screwball@hell.org]</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:t xml:space="preserve">TO: Mugwort</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">&lt;muggy_as@hell.org&gt;</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:t xml:space="preserve">SUBJECT: RE:RE: xxxxxxxxxxxxx</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">You know, these screens really are unpredictable, and we still havnt been able to anticipate all the outcomes from their use with these cattle-like humans. But finally, we''ve caught a big break! These last months I was worried our screen weren''t reaching your patient, but my fears have apparently been unfounded.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">The ultimate goal of these tactics and devices like screens are to bring our patient into various states where he is most likely to stay away from our natural Enemy. Despair is among the most effective of the states that we can lead them into and your man has crossed that threshold brilliantly.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">As you know from your orientation program, Hell.org has legions of departments, many of them siloed off from each other. As part of that new coporate collaboration initiative that none of us were excited about, some of these departments have been teaming up for knowledge transfer, as the consultants called it.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">As a result of this, the climate change department and our screens department have been working tirelessly to expand the environmentalism work that''s been stagnant since the 1970''s. Not only did they up the ante by rebranding environmentalism with the more serious tag, "climate change", they''ve added the term into the current basket of untouchables like Darwinism, moral reletavism, etc.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Now I''ll pre-emptively answer the question you''re sure to ask: yes, Mugwort, even though the Enemy orders his humans to be good stewards to nature, we can and have perverted this command. We cannot create anything; no new feelings, no new sensations. But we can pervert aspects of this world the Enemy has meant for good, to serve our purposes.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Case-in-point: while we should be indifferent or even apprehensive about humans</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">actually</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">taking care of their planet, we should encourage, by all means, the support of this climate change movement. For one thing, when the humans think of saving the world as a whole, they entirely miss the trees over their heads. Keeping them focused on the</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">movement</w:t></w:r><w:r><w:t xml:space="preserve">, rather than nature outside their front door, is among the most hilariously hypocritical things we tempters get to watch. Nature turns into an abstract "thing", an ephemeral idea, rather than a collection of carbon, smells, sounds and life that ultimately points to</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">him</w:t></w:r><w:r><w:t xml:space="preserve">.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Our Enemy put his fingerprints all over nature. A human born on an island who never hears His message would still know Him from his handiword all around. The fibonacci spiral in the conch shell, the movement and pace of the tides, the stars in the sky and the seasons for everything under the life-giving sun, perfectly placed to sustain life and illuminate His face.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">The human enthralled in stopping climate change doesn''t notice this because nature has been reduced to a boogyman -- something that they angered enough to kill them all.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">The point is, do what you can to prevent your patient from going on walks. Whenever he sees a beautiful sunset, redirect him to seeth over the thoughtless oil companies that are one spill away from destroying it. Don''t render him able to enjoy these things; only to think of these wonderful things disappearing.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">You''ve noted that your patient has began to take climate change pretty seriously. Wunderbar. The more he watches videos and speeches that confirm his position, the more entrenched he''ll become. You also noted that he has become more anxious as of late. To be sure, this is most likely from his interest in the topics at-hand.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">The beautiful thing about the climate change angle is the scope of the problem seems huge, impossible to fathom, so it leads these fools into Despair, which is the state we want to keep them in indefinitely.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">When a human is flooded with Despair, they display incredible distrust in the Enemy. While in this state, they are incapable of feeling joy, that most disgusting of human emotions. For some reason, the Enemy wishes these beings to experience joy every day of their petty, short lives. You can see how a cancerous mess of despair can help us in our efforts tremendously.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Keep up the pressure in this area, and you''ll make something of yourself yet.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Screwtape | Executive Tempter</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:t xml:space="preserve">Hell.org</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t xml:space="preserve">Chapter 5</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/></w:pPr></w:p>'
$cursor.InsertXML($xmlSnippet)

# --- Re-locate the two "Chapter 5" headings now present: the first is the OLD heading
#     (about to be renamed to "Chapter 4"), the second is the NEW heading we just inserted. ---
$findRng2 = $d.Content
$findRng2.Find.Execute("Chapter 5") | Out-Null
$firstHeadingPara = $findRng2.Paragraphs(1)
$firstHeadingRange = $d.Range($firstHeadingPara.Range.Start, $firstHeadingPara.Range.Start)

$findRng2.Collapse(0)
$findRng2.Find.Execute("Chapter 5") | Out-Null
$secondHeadingPara = $findRng2.Paragraphs(1)
$secondHeadingRange = $d.Range($secondHeadingPara.Range.Start, $secondHeadingPara.Range.Start)

# --- Remove the trailing placeholder paragraph that sits right after the new heading. ---
$placeholderPara = $secondHeadingPara.Next()
$placeholderPara.Range.Delete()

# --- Move the "chapter-5" bookmark onto the new heading, then add a fresh "chapter-4"
#     bookmark where the old heading used to be (Bookmarks.Add relocates a bookmark whose
#     name already exists rather than duplicating it). ---
$d.Bookmarks.Add("chapter-5", $secondHeadingRange)
$d.Bookmarks.Add("chapter-4", $firstHeadingRange)

# --- Rename the original heading text from "Chapter 5" to "Chapter 4". ---
$firstHeadingPara.Range.Text = "Chapter 4"

# --- Apply the VerbatimChar character style to the Mugwort email address run that InsertXML
#     could not carry over (it drops <w:rStyle> references on insertion). ---
$styleRng = $d.Content
$styleRng.Find.Execute("<muggy_as@hell.org>") | Out-Null
$styleRng.Style = "VerbatimChar"

Write-Output "Chapter 4 inserted and Chapter 5 relocated successfully."
